$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting Mutual Fund..QoQ right by one
$ws.Columns("C").Insert()

# Set the header for the new Industry column
$ws.Cells.Item(1, 3).Value = "Industry"

# Fill in Industry values for each data row (rows 2-22)
$ws.Cells.Item(2, 3).Value = "Banks"
$ws.Cells.Item(3, 3).Value = "Construction"
$ws.Cells.Item(4, 3).Value = "Power"
$ws.Cells.Item(5, 3).Value = "Automobiles"
$ws.Cells.Item(6, 3).Value = "Banks"
$ws.Cells.Item(7, 3).Value = "Power"
$ws.Cells.Item(8, 3).Value = "Insurance"
$ws.Cells.Item(9, 3).Value = "Metals & Minerals Trading"
$ws.Cells.Item(10, 3).Value = "Auto Components"
$ws.Cells.Item(11, 3).Value = "Pharmaceuticals & Biotechnology"
$ws.Cells.Item(12, 3).Value = "Finance"
$ws.Cells.Item(13, 3).Value = "Capital Markets"
$ws.Cells.Item(14, 3).Value = "Construction"
$ws.Cells.Item(15, 3).Value = "Finance"
$ws.Cells.Item(16, 3).Value = "Petroleum Products"
$ws.Cells.Item(17, 3).Value = "IT - Software"
$ws.Cells.Item(18, 3).Value = "Food Products"
$ws.Cells.Item(19, 3).Value = "N.A."
$ws.Cells.Item(20, 3).Value = "Insurance"
$ws.Cells.Item(21, 3).Value = "Banks"
$ws.Cells.Item(22, 3).Value = "Power"
